$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = 37
$ws.Range("C74").Value = 1
$ws.Range("D74").Value = "Token has expired"
$ws.Range("A75").Value = 74
$ws.Range("B75").Value = 37
$ws.Range("C75").Value = 2
$ws.Range("D75").Value = "Tokenın süresi doldu"
$ws.Range("D65").Select() | Out-Null
